$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Principal Component Analysis" (index 6): selection changes only
# ------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A1:XFD1").Select()

# ------------------------------------------------------------------
# Sheet "Change Column Type" (index 7): add Orange-widget panel (F:K)
# mirroring the existing Python panel (B:D), matching the layout used
# on the "Principal Component Analysis" sheet.
# ------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

# --- copy cell formatting (styles) from the already-laid-out sheet6 ---
$ws6.Range("C1").Copy()
$ws7.Range("C1").PasteSpecial(-4122)
$ws7.Range("G1").PasteSpecial(-4122)
$ws7.Range("K1").PasteSpecial(-4122)

$ws6.Range("D1").Copy()
$ws7.Range("D1").PasteSpecial(-4122)
$ws7.Range("H1").PasteSpecial(-4122)

$ws6.Range("B2").Copy()
$ws7.Range("B2:C2").PasteSpecial(-4122)
$ws7.Range("B7:C7").PasteSpecial(-4122)

$ws6.Range("D2").Copy()
$ws7.Range("D2").PasteSpecial(-4122)

$ws6.Range("F3:G3").Copy()
$ws7.Range("F2:G2").PasteSpecial(-4122)
$ws7.Range("F3:G5").PasteSpecial(-4122)

$ws6.Range("H3").Copy()
$ws7.Range("H2:H5").PasteSpecial(-4122)

$ws6.Range("B3:C3").Copy()
$ws7.Range("B3:C6").PasteSpecial(-4122)

$ws6.Range("D3").Copy()
$ws7.Range("D3:D6").PasteSpecial(-4122)

$ws6.Range("F2:G2").Copy()
$ws7.Range("F6:G6").PasteSpecial(-4122)

$ws6.Range("H6").Copy()
$ws7.Range("H6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- values ---
$ws7.Range("C1").Value = "Python"
$ws7.Range("G1").Value = "Orange"
$ws7.Range("K1").Value = "Data Polish"

$ws7.Range("B2").Value = "Action"
$ws7.Range("C2").Value = "Time"
$ws7.Range("D2").Value = "Content"
$ws7.Range("F2").Value = "Action"
$ws7.Range("G2").Value = "Time"
$ws7.Range("H2").Value = "Content"

$ws7.Range("B3").Value = "Upload CSV"
$ws7.Range("C3").Value = "5 min"
$ws7.Range("D3").Value = "df = pd.read_csv('file.csv')"
$ws7.Range("F3").Value = "Load Data"
$ws7.Range("G3").Value = "2 min"
$ws7.Range("H3").Value = "Use 'File' widget to load the dataset"

$ws7.Range("B4").Value = "Identify Column & Type"
$ws7.Range("C4").Value = "2 min"
$ws7.Range("D4").Value = "Review the data with df.dtypes"
$ws7.Range("F4").Value = "Edit Domain"
$ws7.Range("G4").Value = "2 min"
$ws7.Range("H4").Value = "Use 'Edit Domain' widget to rename columns"

$ws7.Range("B5").Value = "Convert Data Type"
$ws7.Range("C5").Value = "2 min"
$ws7.Range("D5").Value = "df['column'] = df['column'].astype('desired_type')"
$ws7.Range("F5").Value = "Verify Changes"
$ws7.Range("G5").Value = "1 min"
$ws7.Range("H5").Value = "Inspect with 'Data Table' widget"

$ws7.Range("B6").Value = "Verify Changes"
$ws7.Range("C6").Value = "1 min"
$ws7.Range("D6").Value = "df.dtypes to confirm the change"
$ws7.Range("F6").Value = "Overall"
$ws7.Range("G6").Value = "5 min"

$ws7.Range("B7").Value = "Overall"
$ws7.Range("C7").Value = "10 min"

# --- row heights ---
$ws7.Rows.Item(1).RowHeight = 19
$ws7.Rows.Item(2).RowHeight = 18
$ws7.Rows.Item(3).RowHeight = 72
$ws7.Rows.Item(4).RowHeight = 90
$ws7.Rows.Item(5).RowHeight = 108
$ws7.Rows.Item(6).RowHeight = 54
$ws7.Rows.Item(7).RowHeight = 17

# --- column widths matching the D/H "wrap" columns elsewhere in the workbook ---
$ws7.Columns.Item(4).ColumnWidth = 10
$ws7.Columns.Item(8).ColumnWidth = 10

# --- view / selection ---
$ws7.Range("A1:XFD1").Select()

# ------------------------------------------------------------------
# Sheet "Rename Column" (index 8): add Orange-widget panel (F:K)
# ------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)

# --- copy cell formatting (styles) from sheet6 ---
$ws6.Range("C1").Copy()
$ws8.Range("C1").PasteSpecial(-4122)
$ws8.Range("G1").PasteSpecial(-4122)
$ws8.Range("K1").PasteSpecial(-4122)

$ws6.Range("D1").Copy()
$ws8.Range("D1").PasteSpecial(-4122)
$ws8.Range("H1").PasteSpecial(-4122)

$ws6.Range("B2").Copy()
$ws8.Range("B2:C2").PasteSpecial(-4122)
$ws8.Range("B6:C6").PasteSpecial(-4122)

$ws6.Range("D2").Copy()
$ws8.Range("D2").PasteSpecial(-4122)

$ws6.Range("F2:G2").Copy()
$ws8.Range("F2:G2").PasteSpecial(-4122)
$ws8.Range("F6:G6").PasteSpecial(-4122)

$ws6.Range("H2").Copy()
$ws8.Range("H2").PasteSpecial(-4122)

$ws6.Range("B3:C3").Copy()
$ws8.Range("B3:C5").PasteSpecial(-4122)

$ws6.Range("D3").Copy()
$ws8.Range("D3:D5").PasteSpecial(-4122)

$ws6.Range("F3:G3").Copy()
$ws8.Range("F3:G5").PasteSpecial(-4122)

$ws6.Range("H3").Copy()
$ws8.Range("H3:H5").PasteSpecial(-4122)

$ws6.Range("H6").Copy()
$ws8.Range("H6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- values ---
$ws8.Range("C1").Value = "Python"
$ws8.Range("G1").Value = "Orange"
$ws8.Range("K1").Value = "Data Polish"

$ws8.Range("B2").Value = "Action"
$ws8.Range("C2").Value = "Time"
$ws8.Range("D2").Value = "Content"
$ws8.Range("F2").Value = "Action"
$ws8.Range("G2").Value = "Time"
$ws8.Range("H2").Value = "Content"

$ws8.Range("B3").Value = "Upload CSV"
$ws8.Range("C3").Value = "5 min"
$ws8.Range("D3").Value = "df = pd.read_csv('file.csv')"
$ws8.Range("F3").Value = "Load Data"
$ws8.Range("G3").Value = "2 min"
$ws8.Range("H3").Value = "Use 'File' widget to load the dataset"

$ws8.Range("B4").Value = "Rename Column"
$ws8.Range("C4").Value = "1 min"
$ws8.Range("D4").Value = "df.rename(columns={'old_name': 'new_name'}, inplace=True)"
$ws8.Range("F4").Value = "Edit Domain"
$ws8.Range("G4").Value = "2 min"
$ws8.Range("H4").Value = "Use 'Edit Domain' widget to rename columns"

$ws8.Range("B5").Value = "Verify Changes"
$ws8.Range("C5").Value = "1 min"
$ws8.Range("D5").Value = "df.head() to check new column names"
$ws8.Range("F5").Value = "Verify Changes"
$ws8.Range("G5").Value = "1 min"
$ws8.Range("H5").Value = "Inspect with 'Data Table' widget"

$ws8.Range("B6").Value = "Overall"
$ws8.Range("C6").Value = "7 min"
$ws8.Range("F6").Value = "Overall"
$ws8.Range("G6").Value = "5 min"

# --- row heights ---
$ws8.Rows.Item(1).RowHeight = 19
$ws8.Rows.Item(2).RowHeight = 18
$ws8.Rows.Item(3).RowHeight = 72
$ws8.Rows.Item(4).RowHeight = 126
$ws8.Rows.Item(5).RowHeight = 90
$ws8.Rows.Item(6).RowHeight = 17

# --- column width matching the D "wrap" column elsewhere in the workbook ---
$ws8.Columns.Item(4).ColumnWidth = 10

# --- view / selection ---
$ws8.Range("F2:H2").Select()

# ------------------------------------------------------------------
# Finally, make "Rename Column" the active tab (matches workbookView
# activeTab) and scroll the tab strip so it (and a few sheets before
# it) are visible.
# ------------------------------------------------------------------
$ws8.Activate()
